$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Delete the oldest reporting-period column (D) entirely.
#    This shifts all later columns (old E..M) one slot to the left (new D..L)
#    and keeps their row-by-row styling/widths intact.
$ws.Range("D:D").EntireColumn.Delete()

# 2. The publish-date cell that used to sit in column J (now shifted to column I)
#    is being re-labelled for the newly issued report revision.
$ws.Range("I9").Value = "1402-02-27 (7)"

# 3. Append the newest reporting period as the new last column (M).
# First, match the formatting of the new M cells to their left neighbor (L), which
# already carries the correct per-row formatting (header/date/data borders+fills).
$ws.Range("L1:L28").Copy()
$ws.Range("M1:M28").PasteSpecial(-4122)
$ws.Columns("M").ColumnWidth = $ws.Columns("L").ColumnWidth

# Now fill in the values on top of that formatting.
$ws.Range("M8").Value = "12 ماهه منتهی به 1401/12"

# Force text so the pure yyyy-mm-dd-looking date string isn't auto-parsed into a date serial,
# then re-paste L9's formats so the temporary text numberformat doesn't stick (keeps border/fill/font
# and General numFmtId identical to the rest of row 9, matching the other (non-pure-date) date cells).
$ws.Range("M9").NumberFormat = "@"
$ws.Range("M9").Value = "1402-02-27"
$ws.Range("L9").Copy()
$ws.Range("M9").PasteSpecial(-4122)

# 4. Fill in the new column M with the income-statement figures for the new period.
$ws.Range("M11").Value = 8500
$ws.Range("M12").Value = -6766
$ws.Range("M13").Value = 1734
$ws.Range("M14").Value = -270
$ws.Range("M15").Value = "-"
$ws.Range("M16").Value = 97
$ws.Range("M17").Value = 1561
$ws.Range("M18").Value = -745
$ws.Range("M19").Value = 83
$ws.Range("M20").Value = 899
$ws.Range("M21").Value = -113
$ws.Range("M22").Value = 785
$ws.Range("M23").Value = "-"
$ws.Range("M24").Value = 785
$ws.Range("M25").Value = 0
$ws.Range("M26").Value = 2196
$ws.Range("M27").Value = 0
